# Auto-generated edit script: apply numeric corrections to profit-calculation
# columns (H, I, J, K, L, M, N) across multiple craft-recipe worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 10916.5
$ws.Range("I40").Value = 8099.8
$ws.Range("J40").Value = 25000
$ws.Range("K40").Value = 8099.8
$ws.Range("L40").Value = 25000
$ws.Range("M40").Value = -7924.8
$ws.Range("N40").Value = -25350
$ws.Range("H86").Value = 5179.95
$ws.Range("I86").Value = 5198.7144
$ws.Range("K86").Value = 5198.7144
$ws.Range("M86").Value = -4075.7144
$ws.Range("H88").Value = 7830.6665
$ws.Range("J88").Value = 7830.6665
$ws.Range("L88").Value = 7830.6665
$ws.Range("N88").Value = -8642.666499999999
$ws.Range("H89").Value = 5179.95
$ws.Range("I89").Value = 5198.7144
$ws.Range("K89").Value = 25993.572
$ws.Range("M89").Value = -20377.572
$ws.Range("H91").Value = 7830.6665
$ws.Range("J91").Value = 7830.6665
$ws.Range("L91").Value = 7830.6665
$ws.Range("N91").Value = -10638.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4079.5789
$ws.Range("I61").Value = 3768.9285
$ws.Range("K61").Value = 3768.9285
$ws.Range("M61").Value = -3556.9285
$ws.Range("H74").Value = 2571
$ws.Range("I74").Value = 2599.25
$ws.Range("K74").Value = 2599.25
$ws.Range("M74").Value = -1725.25
$ws.Range("H77").Value = 2571
$ws.Range("I77").Value = 2599.25
$ws.Range("K77").Value = 12996.25
$ws.Range("M77").Value = -8628.25
$ws.Range("H132").Value = 6858.96
$ws.Range("I132").Value = 6703.773
$ws.Range("K132").Value = 20111.319
$ws.Range("M132").Value = -17581.319
$ws.Range("H134").Value = 80999.5
$ws.Range("J134").Value = 80999.5
$ws.Range("L134").Value = 80999.5
$ws.Range("N134").Value = -91139.5
$ws.Range("H135").Value = 78999
$ws.Range("J135").Value = 78999
$ws.Range("L135").Value = 78999
$ws.Range("N135").Value = -89139
$ws.Range("H136").Value = 4079.5789
$ws.Range("I136").Value = 3768.9285
$ws.Range("K136").Value = 11306.7855
$ws.Range("M136").Value = -8756.7855
$ws.Range("H137").Value = 80000
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2415.1
$ws.Range("I86").Value = 1225.0769
$ws.Range("J86").Value = 4625.143
$ws.Range("K86").Value = 1225.0769
$ws.Range("L86").Value = 4625.143
$ws.Range("M86").Value = -102.0769
$ws.Range("N86").Value = -6871.143
$ws.Range("H89").Value = 2415.1
$ws.Range("I89").Value = 1225.0769
$ws.Range("J89").Value = 4625.143
$ws.Range("K89").Value = 6125.3845
$ws.Range("L89").Value = 23125.715
$ws.Range("M89").Value = -509.3845000000001
$ws.Range("N89").Value = -34357.715
$ws.Range("H134").Value = 2906.625
$ws.Range("I134").Value = 2904.639
$ws.Range("K134").Value = 8713.917000000001
$ws.Range("M134").Value = -6178.917000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1094.4
$ws.Range("I4").Value = 740.8
$ws.Range("J4").Value = 1801.6
$ws.Range("K4").Value = 740.8
$ws.Range("L4").Value = 1801.6
$ws.Range("M4").Value = -628.8
$ws.Range("N4").Value = -2025.6
$ws.Range("H58").Value = 4219.483
$ws.Range("J58").Value = 4028.1667
$ws.Range("L58").Value = 4028.1667
$ws.Range("N58").Value = -4434.1667
$ws.Range("H122").Value = 4371.5884
$ws.Range("I122").Value = 4329.7856
$ws.Range("J122").Value = 4566.6665
$ws.Range("K122").Value = 12989.3568
$ws.Range("L122").Value = 13699.9995
$ws.Range("M122").Value = -10539.3568
$ws.Range("N122").Value = -18599.9995
$ws.Range("H136").Value = 4219.483
$ws.Range("J136").Value = 4028.1667
$ws.Range("L136").Value = 12084.5001
$ws.Range("N136").Value = -17184.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1251.7
$ws.Range("I17").Value = 146
$ws.Range("J17").Value = 2910.25
$ws.Range("K17").Value = 438
$ws.Range("L17").Value = 8730.75
$ws.Range("M17").Value = -269
$ws.Range("N17").Value = -9068.75
$ws.Range("H29").Value = 10938.3
$ws.Range("J29").Value = 983.8570999999999
$ws.Range("L29").Value = 2951.5713
$ws.Range("N29").Value = -3505.5713
$ws.Range("H61").Value = 437.6
$ws.Range("I61").Value = 422
$ws.Range("K61").Value = 1266
$ws.Range("M61").Value = -1051
$ws.Range("H62").Value = 9248.5
$ws.Range("J62").Value = 10999
$ws.Range("L62").Value = 32997
$ws.Range("N62").Value = -34369
$ws.Range("H65").Value = 9248.5
$ws.Range("J65").Value = 10999
$ws.Range("N65").Value = -105855
$ws.Range("H119").Value = 6344
$ws.Range("I119").Value = 2930.375
$ws.Range("K119").Value = 8791.125
$ws.Range("M119").Value = -3953.125
$ws.Range("H123").Value = 4997
$ws.Range("J123").Value = 4400
$ws.Range("L123").Value = 13200
$ws.Range("N123").Value = -18100

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8898.857
$ws.Range("I70").Value = 8823.75
$ws.Range("K70").Value = 8823.75
$ws.Range("M70").Value = -8553.75
$ws.Range("H73").Value = 8898.857
$ws.Range("I73").Value = 8823.75
$ws.Range("K73").Value = 8823.75
$ws.Range("M73").Value = -7887.75
$ws.Range("H132").Value = 4682.3706
$ws.Range("I132").Value = 5274.1113
$ws.Range("J132").Value = 3498.889
$ws.Range("K132").Value = 15822.3339
$ws.Range("L132").Value = 10496.667
$ws.Range("M132").Value = -13292.3339
$ws.Range("N132").Value = -15556.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4970.2666
$ws.Range("I46").Value = 5839
$ws.Range("J46").Value = 2581.25
$ws.Range("K46").Value = 5839
$ws.Range("L46").Value = 2581.25
$ws.Range("M46").Value = -5651
$ws.Range("N46").Value = -2957.25
$ws.Range("H136").Value = 1778
$ws.Range("I136").Value = 1773.8334
$ws.Range("J136").Value = 1786.3334
$ws.Range("K136").Value = 5321.5002
$ws.Range("L136").Value = 5359.0002
$ws.Range("M136").Value = -2771.5002
$ws.Range("N136").Value = -10459.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 6023.2
$ws.Range("I9").Value = 4279
$ws.Range("K9").Value = 4279
$ws.Range("M9").Value = -4139
$ws.Range("H132").Value = 3116.4075
$ws.Range("I132").Value = 3395.348
$ws.Range("K132").Value = 10186.044
$ws.Range("M132").Value = -7656.044
$ws.Range("H136").Value = 3757.9285
$ws.Range("I136").Value = 3606
$ws.Range("K136").Value = 10818
$ws.Range("M136").Value = -8268
